$d = $word.ActiveDocument

# Update the date paragraph (first paragraph, centered, Arial)
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.End = $dateRange.End - 1
$dateRange.Text = "2024-10-12 Saturday"

# Update the table cells (20 rows x 5 columns), in row-major order
$t = $d.Tables.Item(1)
$values = @(
    "42+32=74", "14+17=31", "57-54=3", "86-14=72", "50-46=4",
    "23+32=55", "31+7=38", "66-53=13", "47-9=38", "94-11=83",
    "87-72=15", "76-23=53", "4+19=23", "89-24=65", "13+85=98",
    "54-13=41", "47+8=55", "46+34=80", "64-33=31", "76-5=71",
    "23+74=97", "81-27=54", "81-8=73", "30+42=72", "50-11=39",
    "50-15=35", "38+38=76", "96-3=93", "53-1=52", "20+53=73",
    "81-62=19", "95-59=36", "66-3=63", "82-25=57", "81+2=83",
    "82-82=0", "12+57=69", "42-35=7", "82+15=97", "26-0=26",
    "50-49=1", "44+46=90", "24+22=46", "52-2=50", "22+2=24",
    "11+72=83", "82-35=47", "50+35=85", "4+75=79", "63-60=3",
    "10+55=65", "62-2=60", "30-1=29", "40+53=93", "53+23=76",
    "83-66=17", "2+43=45", "77+0=77", "51+40=91", "64+9=73",
    "16+82=98", "89-64=25", "53-45=8", "83-47=36", "84-30=54",
    "52-21=31", "57-23=34", "68-54=14", "37+31=68", "70-15=55",
    "66+20=86", "19+41=60", "5+68=73", "13+51=64", "78-23=55",
    "34+58=92", "23+12=35", "37-1=36", "22+39=61", "10+15=25",
    "77-1=76", "36+5=41", "81-76=5", "5+61=66", "10+87=97",
    "38-20=18", "19+9=28", "24-13=11", "18-12=6", "34-0=34",
    "57+40=97", "71-68=3", "8+48=56", "40+27=67", "80+9=89",
    "84-76=8", "32+8=40", "94-76=18", "77-16=61", "43+37=80"
)

$idx = 0
for ($row = 1; $row -le 20; $row++) {
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $r = $cell.Range
        $r.End = $r.End - 1
        $r.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Done. idx=" $idx